# Auto-generated schedule update script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'Monday'
$ws.Cells.Item(2, 2).Value = '08:00-11:00'
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 'Alexander, Alejandro L, Adi'
$ws.Cells.Item(2, 5).Value = 'Alexander, Adi'

$ws.Cells.Item(3, 1).Value = 'Monday'
$ws.Cells.Item(3, 2).Value = '11:00-14:00'
$ws.Cells.Item(3, 3).Value = 1
$ws.Cells.Item(3, 4).Value = 'Jaxon, Ben Kairouz, Kamsi'
$ws.Cells.Item(3, 5).Value = 'Jaxon, Ben Kairouz'

$ws.Cells.Item(4, 1).Value = 'Monday'
$ws.Cells.Item(4, 2).Value = '14:00-17:00'
$ws.Cells.Item(4, 3).Value = 2
$ws.Cells.Item(4, 4).Value = 'Henry, Jamari Pitchford, Harry Corbin'
$ws.Cells.Item(4, 5).Value = 'Jamari Pitchford, Harry Corbin'

$ws.Cells.Item(5, 1).Value = 'Monday'
$ws.Cells.Item(5, 2).Value = '17:00-20:00'
$ws.Cells.Item(5, 3).Value = 3
$ws.Cells.Item(5, 4).Value = 'Thor Waguespack, Ali Awada, Alejandro E. Ulvert'
$ws.Cells.Item(5, 5).Value = 'Thor Waguespack, Alejandro E. Ulvert'

$ws.Cells.Item(6, 1).Value = 'Monday'
$ws.Cells.Item(6, 2).Value = '20:00-23:00'
$ws.Cells.Item(6, 3).Value = 4
$ws.Cells.Item(6, 4).Value = ''
$ws.Cells.Item(6, 5).Value = ''

$ws.Cells.Item(7, 1).Value = 'Monday'
$ws.Cells.Item(7, 2).Value = '23:00-02:00'
$ws.Cells.Item(7, 3).Value = 5
$ws.Cells.Item(7, 4).Value = ''
$ws.Cells.Item(7, 5).Value = ''

$ws.Cells.Item(8, 1).Value = 'Monday'
$ws.Cells.Item(8, 2).Value = '02:00-08:00 (Night)'
$ws.Cells.Item(8, 3).Value = 6
$ws.Cells.Item(8, 4).Value = 'Ezana, Paul'
$ws.Cells.Item(8, 5).Value = 'Paul'

$ws.Cells.Item(9, 1).Value = 'Tuesday'
$ws.Cells.Item(9, 2).Value = '08:00-11:00'
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(9, 4).Value = 'Paul, Alexander, George Ryckman'
$ws.Cells.Item(9, 5).Value = 'Paul, Alexander, George Ryckman'

$ws.Cells.Item(10, 1).Value = 'Tuesday'
$ws.Cells.Item(10, 2).Value = '11:00-14:00'
$ws.Cells.Item(10, 3).Value = 1
$ws.Cells.Item(10, 4).Value = 'Alejandro Espinosa, Jamari Pitchford, Blake Steel'
$ws.Cells.Item(10, 5).Value = 'Alejandro Espinosa, Jamari Pitchford, Blake Steel'

$ws.Cells.Item(11, 1).Value = 'Tuesday'
$ws.Cells.Item(11, 2).Value = '14:00-17:00'
$ws.Cells.Item(11, 3).Value = 2
$ws.Cells.Item(11, 4).Value = 'Josh Greene, Matheo, Jake Dieterich'
$ws.Cells.Item(11, 5).Value = 'Josh Greene, Matheo'

$ws.Cells.Item(12, 1).Value = 'Tuesday'
$ws.Cells.Item(12, 2).Value = '17:00-20:00'
$ws.Cells.Item(12, 3).Value = 3
$ws.Cells.Item(12, 4).Value = 'Jaxon, Gabe Heller, Alejandro E. Ulvert'
$ws.Cells.Item(12, 5).Value = 'Jaxon, Gabe Heller, Alejandro E. Ulvert'

$ws.Cells.Item(13, 1).Value = 'Tuesday'
$ws.Cells.Item(13, 2).Value = '20:00-23:00'
$ws.Cells.Item(13, 3).Value = 4
$ws.Cells.Item(13, 4).Value = ''
$ws.Cells.Item(13, 5).Value = ''

$ws.Cells.Item(14, 1).Value = 'Tuesday'
$ws.Cells.Item(14, 2).Value = '23:00-02:00'
$ws.Cells.Item(14, 3).Value = 5
$ws.Cells.Item(14, 4).Value = ''
$ws.Cells.Item(14, 5).Value = ''

$ws.Cells.Item(15, 1).Value = 'Tuesday'
$ws.Cells.Item(15, 2).Value = '02:00-08:00 (Night)'
$ws.Cells.Item(15, 3).Value = 6
$ws.Cells.Item(15, 4).Value = 'Jack Mogelof, Henry'
$ws.Cells.Item(15, 5).Value = 'Jack Mogelof'

$ws.Cells.Item(16, 1).Value = 'Wednesday'
$ws.Cells.Item(16, 2).Value = '08:00-11:00'
$ws.Cells.Item(16, 3).Value = 0
$ws.Cells.Item(16, 4).Value = 'Ezana, Josh Greene, Thor Waguespack'
$ws.Cells.Item(16, 5).Value = 'Josh Greene, Thor Waguespack'

$ws.Cells.Item(17, 1).Value = 'Wednesday'
$ws.Cells.Item(17, 2).Value = '11:00-14:00'
$ws.Cells.Item(17, 3).Value = 1
$ws.Cells.Item(17, 4).Value = 'Gabe Heller, Jake Dieterich, Noah Yaffe'
$ws.Cells.Item(17, 5).Value = 'Gabe Heller, Noah Yaffe'

$ws.Cells.Item(18, 1).Value = 'Wednesday'
$ws.Cells.Item(18, 2).Value = '14:00-17:00'
$ws.Cells.Item(18, 3).Value = 2
$ws.Cells.Item(18, 4).Value = 'Henry, Alejandro Espinosa, Edu'
$ws.Cells.Item(18, 5).Value = 'Alejandro Espinosa, Edu'

$ws.Cells.Item(19, 1).Value = 'Wednesday'
$ws.Cells.Item(19, 2).Value = '17:00-20:00'
$ws.Cells.Item(19, 3).Value = 3
$ws.Cells.Item(19, 4).Value = 'Ben Kairouz, Kamsi, Harry Corbin'
$ws.Cells.Item(19, 5).Value = 'Ben Kairouz, Harry Corbin'

$ws.Cells.Item(20, 1).Value = 'Wednesday'
$ws.Cells.Item(20, 2).Value = '20:00-23:00'
$ws.Cells.Item(20, 3).Value = 4
$ws.Cells.Item(20, 4).Value = ''
$ws.Cells.Item(20, 5).Value = ''

$ws.Cells.Item(21, 1).Value = 'Wednesday'
$ws.Cells.Item(21, 2).Value = '23:00-02:00'
$ws.Cells.Item(21, 3).Value = 5
$ws.Cells.Item(21, 4).Value = ''
$ws.Cells.Item(21, 5).Value = ''

$ws.Cells.Item(22, 1).Value = 'Wednesday'
$ws.Cells.Item(22, 2).Value = '02:00-08:00 (Night)'
$ws.Cells.Item(22, 3).Value = 6
$ws.Cells.Item(22, 4).Value = 'Ali Awada, Blake Steel'
$ws.Cells.Item(22, 5).Value = 'Blake Steel'

$ws.Cells.Item(23, 1).Value = 'Thursday'
$ws.Cells.Item(23, 2).Value = '08:00-11:00'
$ws.Cells.Item(23, 3).Value = 0
$ws.Cells.Item(23, 4).Value = 'Jack Mogelof, George Ryckman, Noah Yaffe'
$ws.Cells.Item(23, 5).Value = 'Jack Mogelof, George Ryckman, Noah Yaffe'

$ws.Cells.Item(24, 1).Value = 'Thursday'
$ws.Cells.Item(24, 2).Value = '11:00-14:00'
$ws.Cells.Item(24, 3).Value = 1
$ws.Cells.Item(24, 4).Value = 'Ezana, Josh Greene, Thor Waguespack'
$ws.Cells.Item(24, 5).Value = 'Josh Greene, Thor Waguespack'

$ws.Cells.Item(25, 1).Value = 'Thursday'
$ws.Cells.Item(25, 2).Value = '14:00-17:00'
$ws.Cells.Item(25, 3).Value = 2
$ws.Cells.Item(25, 4).Value = 'Gabe Heller, Ali Awada, Alejandro E. Ulvert'
$ws.Cells.Item(25, 5).Value = 'Gabe Heller, Alejandro E. Ulvert'

$ws.Cells.Item(26, 1).Value = 'Thursday'
$ws.Cells.Item(26, 2).Value = '17:00-20:00'
$ws.Cells.Item(26, 3).Value = 3
$ws.Cells.Item(26, 4).Value = 'Matheo, Edu, Adi'
$ws.Cells.Item(26, 5).Value = 'Matheo, Edu, Adi'

$ws.Cells.Item(27, 1).Value = 'Thursday'
$ws.Cells.Item(27, 2).Value = '20:00-23:00'
$ws.Cells.Item(27, 3).Value = 4
$ws.Cells.Item(27, 4).Value = ''
$ws.Cells.Item(27, 5).Value = ''

$ws.Cells.Item(28, 1).Value = 'Thursday'
$ws.Cells.Item(28, 2).Value = '23:00-02:00'
$ws.Cells.Item(28, 3).Value = 5
$ws.Cells.Item(28, 4).Value = ''
$ws.Cells.Item(28, 5).Value = ''

$ws.Cells.Item(29, 1).Value = 'Thursday'
$ws.Cells.Item(29, 2).Value = '02:00-08:00 (Night)'
$ws.Cells.Item(29, 3).Value = 6
$ws.Cells.Item(29, 4).Value = 'Jaxon, Alejandro L'
$ws.Cells.Item(29, 5).Value = 'Jaxon'

$ws.Cells.Item(30, 1).Value = 'Friday'
$ws.Cells.Item(30, 2).Value = '08:00-11:00'
$ws.Cells.Item(30, 3).Value = 0
$ws.Cells.Item(30, 4).Value = 'Ben Kairouz, Alexander, Blake Steel'
$ws.Cells.Item(30, 5).Value = 'Ben Kairouz, Alexander, Blake Steel'

$ws.Cells.Item(31, 1).Value = 'Friday'
$ws.Cells.Item(31, 2).Value = '11:00-14:00'
$ws.Cells.Item(31, 3).Value = 1
$ws.Cells.Item(31, 4).Value = 'Jamari Pitchford, George Ryckman, Jake Dieterich'
$ws.Cells.Item(31, 5).Value = 'Jamari Pitchford, George Ryckman'

$ws.Cells.Item(32, 1).Value = 'Friday'
$ws.Cells.Item(32, 2).Value = '14:00-17:00'
$ws.Cells.Item(32, 3).Value = 2
$ws.Cells.Item(32, 4).Value = 'Kamsi, Paul, Alejandro Espinosa'
$ws.Cells.Item(32, 5).Value = 'Paul, Alejandro Espinosa'

$ws.Cells.Item(33, 1).Value = 'Friday'
$ws.Cells.Item(33, 2).Value = '17:00-20:00'
$ws.Cells.Item(33, 3).Value = 3
$ws.Cells.Item(33, 4).Value = 'Matheo, Noah Yaffe, Harry Corbin'
$ws.Cells.Item(33, 5).Value = 'Matheo, Noah Yaffe, Harry Corbin'

$ws.Cells.Item(34, 1).Value = 'Friday'
$ws.Cells.Item(34, 2).Value = '20:00-23:00'
$ws.Cells.Item(34, 3).Value = 4
$ws.Cells.Item(34, 4).Value = ''
$ws.Cells.Item(34, 5).Value = ''

$ws.Cells.Item(35, 1).Value = 'Friday'
$ws.Cells.Item(35, 2).Value = '23:00-02:00'
$ws.Cells.Item(35, 3).Value = 5
$ws.Cells.Item(35, 4).Value = ''
$ws.Cells.Item(35, 5).Value = ''

$ws.Cells.Item(36, 1).Value = 'Friday'
$ws.Cells.Item(36, 2).Value = '02:00-08:00 (Night)'
$ws.Cells.Item(36, 3).Value = 6
$ws.Cells.Item(36, 4).Value = 'Ezana, Alejandro E. Ulvert'
$ws.Cells.Item(36, 5).Value = 'Alejandro E. Ulvert'

$ws.Cells.Item(37, 1).Value = 'Saturday'
$ws.Cells.Item(37, 2).Value = '08:00-11:00'
$ws.Cells.Item(37, 3).Value = 0
$ws.Cells.Item(37, 4).Value = 'Kamsi, Josh Greene, Edu'
$ws.Cells.Item(37, 5).Value = 'Josh Greene, Edu'

$ws.Cells.Item(38, 1).Value = 'Saturday'
$ws.Cells.Item(38, 2).Value = '11:00-14:00'
$ws.Cells.Item(38, 3).Value = 1
$ws.Cells.Item(38, 4).Value = 'Matheo, Alejandro L, Adi'
$ws.Cells.Item(38, 5).Value = 'Matheo, Adi'

$ws.Cells.Item(39, 1).Value = 'Saturday'
$ws.Cells.Item(39, 2).Value = '14:00-17:00'
$ws.Cells.Item(39, 3).Value = 2
$ws.Cells.Item(39, 4).Value = 'Jack Mogelof, Thor Waguespack, Ali Awada'
$ws.Cells.Item(39, 5).Value = 'Jack Mogelof, Thor Waguespack'

$ws.Cells.Item(40, 1).Value = 'Saturday'
$ws.Cells.Item(40, 2).Value = '17:00-20:00'
$ws.Cells.Item(40, 3).Value = 3
$ws.Cells.Item(40, 4).Value = 'Paul, Jamari Pitchford, Blake Steel'
$ws.Cells.Item(40, 5).Value = 'Paul, Jamari Pitchford, Blake Steel'

$ws.Cells.Item(41, 1).Value = 'Saturday'
$ws.Cells.Item(41, 2).Value = '20:00-23:00'
$ws.Cells.Item(41, 3).Value = 4
$ws.Cells.Item(41, 4).Value = ''
$ws.Cells.Item(41, 5).Value = ''

$ws.Cells.Item(42, 1).Value = 'Saturday'
$ws.Cells.Item(42, 2).Value = '23:00-02:00'
$ws.Cells.Item(42, 3).Value = 5
$ws.Cells.Item(42, 4).Value = ''
$ws.Cells.Item(42, 5).Value = ''

$ws.Cells.Item(43, 1).Value = 'Saturday'
$ws.Cells.Item(43, 2).Value = '02:00-08:00 (Night)'
$ws.Cells.Item(43, 3).Value = 6
$ws.Cells.Item(43, 4).Value = 'Gabe Heller, Jake Dieterich'
$ws.Cells.Item(43, 5).Value = 'Gabe Heller'

$ws.Cells.Item(44, 1).Value = 'Sunday'
$ws.Cells.Item(44, 2).Value = '08:00-11:00'
$ws.Cells.Item(44, 3).Value = 0
$ws.Cells.Item(44, 4).Value = 'Jack Mogelof, Henry, George Ryckman'
$ws.Cells.Item(44, 5).Value = 'Jack Mogelof, George Ryckman'

$ws.Cells.Item(45, 1).Value = 'Sunday'
$ws.Cells.Item(45, 2).Value = '11:00-14:00'
$ws.Cells.Item(45, 3).Value = 1
$ws.Cells.Item(45, 4).Value = 'Ben Kairouz, Alejandro L, Adi'
$ws.Cells.Item(45, 5).Value = 'Ben Kairouz, Adi'

$ws.Cells.Item(46, 1).Value = 'Sunday'
$ws.Cells.Item(46, 2).Value = '14:00-17:00'
$ws.Cells.Item(46, 3).Value = 2
$ws.Cells.Item(46, 4).Value = 'Alexander, Noah Yaffe, Harry Corbin'
$ws.Cells.Item(46, 5).Value = 'Alexander, Noah Yaffe, Harry Corbin'

$ws.Cells.Item(47, 1).Value = 'Sunday'
$ws.Cells.Item(47, 2).Value = '17:00-20:00'
$ws.Cells.Item(47, 3).Value = 3
$ws.Cells.Item(47, 4).Value = ''
$ws.Cells.Item(47, 5).Value = ''

$ws.Cells.Item(48, 1).Value = 'Sunday'
$ws.Cells.Item(48, 2).Value = '20:00-23:00'
$ws.Cells.Item(48, 3).Value = 4
$ws.Cells.Item(48, 4).Value = ''
$ws.Cells.Item(48, 5).Value = ''

$ws.Cells.Item(49, 1).Value = 'Sunday'
$ws.Cells.Item(49, 2).Value = '23:00-02:00'
$ws.Cells.Item(49, 3).Value = 5
$ws.Cells.Item(49, 4).Value = ''
$ws.Cells.Item(49, 5).Value = ''

$ws.Cells.Item(50, 1).Value = 'Sunday'
$ws.Cells.Item(50, 2).Value = '02:00-08:00 (Night)'
$ws.Cells.Item(50, 3).Value = 6
$ws.Cells.Item(50, 4).Value = 'Alejandro Espinosa, Edu'
$ws.Cells.Item(50, 5).Value = 'Alejandro Espinosa, Edu'
